$d = $word.ActiveDocument

function Set-ParaXml($paraIndex, [string]$xml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.InsertXML($xml)
}

# Paragraph 28: "A combination of operands and operators make an expression. Even a
# single value is an expression." -> append a new run with the extra sentence.
$p28 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003904C5" w:rsidRDefault="005611BC"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">A </w:t></w:r><w:r w:rsidR="004E4092"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>combination</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of operands and operators make an expression. </w:t></w:r><w:r w:rsidR="004E4092"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Even a single value is an expression.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Expressions executes and gives a value</w:t></w:r></w:p>'
Set-ParaXml 28 $p28

# Paragraph 30: "...between a declaration and an expression?" -> "...between a statement and
# an expression?" split across three runs.
$p30 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003904C5" w:rsidRDefault="004E4092"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">5. In this chapter, assignment statements such as spam = 10 were added. What''s the difference between a </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>statement</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and an expression?</w:t></w:r></w:p>'
Set-ParaXml 30 $p30

# Paragraph 31: Expand the "Answer" paragraph with the new statement/expression explanation.
$p31 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004E4092" w:rsidRDefault="004E4092"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00DC7069"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Answer</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Expression is </w:t></w:r><w:r w:rsidR="00E83A06"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>something which gives a value or a value itself</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>statement</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>something that makes a line, also an expression can be a part of statement or a statement itself</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> In the example 10 is the expression and spam = 10 is the statement.</w:t></w:r></w:p>'
Set-ParaXml 31 $p31

# Paragraph 54/55: move <w:lastRenderedPageBreak/> from the quote run (55) to the
# "10. What is the error..." run (54).
$p54 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003904C5" w:rsidRDefault="004E4092"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>10. What is the error caused by this expression? What would you do about it?</w:t></w:r></w:p>'
Set-ParaXml 54 $p54

$p55 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003904C5" w:rsidRDefault="004E4092"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>''I have eaten '' + 99 + '' burritos.''</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
Set-ParaXml 55 $p55
